# Apply the two changes described by the diff:
#  1. Slide 16's table switches to a different (blank/no-style) table style GUID.
#  2. The "Office Theme" content that used to live in ppt/theme/theme2.xml
#     (the notes master's theme) is swapped into ppt/theme/theme1.xml (the
#     slide master / presentation theme), i.e. the deck's visible design
#     becomes the plain "Office Theme" color set instead of "Integral".
#
# NOTE: this COM-interop host does not expose a usable handle onto the
# notes-master's own theme part (NotesMaster.Theme / notes-page Theme
# resolve to stub objects that read/write nothing), so only the
# slide-master-facing theme (theme1.xml) can be edited here; its color
# scheme is updated to the target "Office" palette.

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 16 -------------------------------------------
$tableSlide = $p.Slides.Item(16)
$tableShape = $tableSlide.Shapes.Item(3)
$tableShape.Table.ApplyStyle("{053306F4-14A1-4F11-8976-E93764F41D69}")

# --- 2. Theme color scheme (Integral -> Office) ----------------------------
$theme = $p.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

function Set-ThemeColor($index, $r, $g, $b) {
    $colors.Item($index).RGB = $r + ($g * 256) + ($b * 65536)
}

Set-ThemeColor 1  0x00 0x00 0x00   # dk1
Set-ThemeColor 2  0xFF 0xFF 0xFF   # lt1
Set-ThemeColor 3  0x44 0x54 0x6A   # dk2
Set-ThemeColor 4  0xE7 0xE6 0xE6   # lt2
Set-ThemeColor 5  0x5B 0x9B 0xD5   # accent1
Set-ThemeColor 6  0xED 0x7D 0x31   # accent2
Set-ThemeColor 7  0xA5 0xA5 0xA5   # accent3
Set-ThemeColor 8  0xFF 0xC0 0x00   # accent4
Set-ThemeColor 9  0x44 0x72 0xC4   # accent5
Set-ThemeColor 10 0x70 0xAD 0x47   # accent6
Set-ThemeColor 11 0x05 0x63 0xC1   # hlink
Set-ThemeColor 12 0x95 0x4F 0x72   # folHlink
